$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("doctors")

# Merge ward codes: "4_2" + "4_1" -> "4_2_4_1" ; "12_2" + "6_3" -> "12_2_6_3"
$ws.Range("B2").Value = "4_2_4_1"
$ws.Range("B3").Value = "4_2_4_1"
$ws.Range("B4").Value = "4_2_4_1"

$ws.Range("B7").Value = "12_2_6_3"
$ws.Range("B8").Value = "12_2_6_3"
$ws.Range("B9").Value = "12_2_6_3"
$ws.Range("B10").Value = "12_2_6_3"

# friday_sunday column flips for rows 4 and 8
$ws.Range("F4").Value = "yes"
$ws.Range("F8").Value = "no"

# Column width adjustments (closest achievable values on Excel's
# character-width pixel grid to the target widths of 9.6 and 13.7)
$ws.Columns.Item(2).ColumnWidth = 8.75
$ws.Columns.Item(3).ColumnWidth = 12.75
